# Fix the "FIRTS_NAME" typo to "FIRST_NAME" in the
# "Nombres y apellidos del titular: ..." placeholder line, and leave the
# text split across three runs (matching the target OOXML) by briefly
# toggling a formatting property on/off over the exact split points —
# this forces the engine to keep the run boundaries without altering the
# visible formatting (all runs keep sz/szCs = 20, no bold).

$d = $word.ActiveDocument

# 1) Correct the typo: FIRTS_NAME -> FIRST_NAME (text-only change).
$null = $d.Content.Find.Execute("FIRTS_NAME", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "FIRST_NAME", 2)

# 2) Re-find the now-corrected line so we can re-split its single merged
#    run into three runs at the exact boundaries the target document uses:
#      "Nombres y apellidos del titular: FIR" | "S" | "T_NAME LAST_NAME"
$line = $d.Content
$line.Find.ClearFormatting()
$found = $line.Find.Execute("Nombres y apellidos del titular: FIRST_NAME LAST_NAME", `
                             $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $lineStart = $line.Start

    # Split boundary #1: right after "...titular: FIR" (37 chars in).
    $splitA = $d.Range($lineStart + 36, $lineStart + 37)
    $splitA.Bold = 1
    $splitA.Bold = 0
}
